# Remove the dependency of the "Dose multiplier interval" columns from the
# pH-meter protocol sheet. In the original layout this value appeared twice
# per dosing stage (column H and column N); both are deleted here, which
# shifts every later column one step to the left for each deletion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the right-most occurrence first so the reference to the left-most
# occurrence ("H") still points at the original column when it is removed.
$ws.Columns("N").EntireColumn.Delete() | Out-Null
$ws.Columns("H").EntireColumn.Delete() | Out-Null

# Restore a sensible selection after the structural edit.
$ws.Range("K8").Select() | Out-Null
